$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells store numbers as text in the source data (runs/balls columns),
# so force Text formatting before writing to avoid Excel auto-converting the
# numeric-looking strings into real numbers.
$cells = @("C2", "D2", "C4", "D4")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("C2").Value = "0"
$ws.Range("D2").Value = "1"
$ws.Range("C4").Value = "1"
$ws.Range("D4").Value = "3"
